$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.017800088035567
$ws.Cells.Item(2, 4).Value = 1.033195546244979
$ws.Cells.Item(2, 5).Value = 1.018894137335163
$ws.Cells.Item(2, 6).Value = 1.026603798390642
$ws.Cells.Item(2, 9).Value = 1.032455306622417
$ws.Cells.Item(2, 10).Value = 1.023012145862818
$ws.Cells.Item(2, 11).Value = 1.035998593965702
$ws.Cells.Item(2, 12).Value = 1.021738998724512
$ws.Cells.Item(2, 13).Value = 1.029425961288676
$ws.Cells.Item(2, 14).Value = 1.024464940451079

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.022068028152732
$ws.Cells.Item(3, 4).Value = 1.033967689226101
$ws.Cells.Item(3, 5).Value = 1.022618182624449
$ws.Cells.Item(3, 6).Value = 1.030636451199123
$ws.Cells.Item(3, 9).Value = 1.032872349603414
$ws.Cells.Item(3, 10).Value = 1.026898778688175
$ws.Cells.Item(3, 11).Value = 1.036579765047623
$ws.Cells.Item(3, 12).Value = 1.025260943140016
$ws.Cells.Item(3, 13).Value = 1.033257457600933
$ws.Cells.Item(3, 14).Value = 1.028357092740849

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.024790381117828
$ws.Cells.Item(4, 4).Value = 1.034459175664764
$ws.Cells.Item(4, 5).Value = 1.024993304244072
$ws.Cells.Item(4, 6).Value = 1.033197377007351
$ws.Cells.Item(4, 9).Value = 1.033133035514955
$ws.Cells.Item(4, 10).Value = 1.029375600737346
$ws.Cells.Item(4, 11).Value = 1.036947043204712
$ws.Cells.Item(4, 12).Value = 1.027505426171509
$ws.Cells.Item(4, 13).Value = 1.035688449974263
$ws.Cells.Item(4, 14).Value = 1.030837432161425

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.025925781571071
$ws.Cells.Item(5, 4).Value = 1.034663882202422
$ws.Cells.Item(5, 5).Value = 1.025983811699965
$ws.Cells.Item(5, 6).Value = 1.034262718673621
$ws.Cells.Item(5, 9).Value = 1.033240466388023
$ws.Cells.Item(5, 10).Value = 1.030408043806647
$ws.Cells.Item(5, 11).Value = 1.037099379581727
$ws.Cells.Item(5, 12).Value = 1.028441036548954
$ws.Cells.Item(5, 13).Value = 1.036699211569928
$ws.Cells.Item(5, 14).Value = 1.031871341418307

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.026115897351467
$ws.Cells.Item(6, 4).Value = 1.034698142048261
$ws.Cells.Item(6, 5).Value = 1.026149661595702
$ws.Cells.Item(6, 6).Value = 1.034440942637748
$ws.Cells.Item(6, 9).Value = 1.033258378695465
$ws.Cells.Item(6, 10).Value = 1.030580887391164
$ws.Cells.Item(6, 11).Value = 1.037124837259005
$ws.Cells.Item(6, 12).Value = 1.028597670106555
$ws.Cells.Item(6, 13).Value = 1.036868273573103
$ws.Cells.Item(6, 14).Value = 1.03204443046054

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.024805587633108
$ws.Cells.Item(7, 4).Value = 1.03446191844057
$ws.Cells.Item(7, 5).Value = 1.025006570492382
$ws.Cells.Item(7, 6).Value = 1.033211656014271
$ws.Cells.Item(7, 9).Value = 1.033134479461276
$ws.Cells.Item(7, 10).Value = 1.029389430524962
$ws.Cells.Item(7, 11).Value = 1.036949086809835
$ws.Cells.Item(7, 12).Value = 1.027517958801846
$ws.Cells.Item(7, 13).Value = 1.035701999520818
$ws.Cells.Item(7, 14).Value = 1.030851281588925

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.019250796406235
$ws.Cells.Item(8, 4).Value = 1.033458202917035
$ws.Cells.Item(8, 5).Value = 1.020160033003057
$ws.Cells.Item(8, 6).Value = 1.027976893379196
$ws.Cells.Item(8, 9).Value = 1.032598167485443
$ws.Cells.Item(8, 10).Value = 1.024333726612973
$ws.Cells.Item(8, 11).Value = 1.036196842267722
$ws.Cells.Item(8, 12).Value = 1.022936560305088
$ws.Cells.Item(8, 13).Value = 1.030731020661109
$ws.Cells.Item(8, 14).Value = 1.025788397997484

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.009146104565222
$ws.Cells.Item(9, 4).Value = 1.031625510542972
$ws.Cells.Item(9, 5).Value = 1.011341455394019
$ws.Cells.Item(9, 6).Value = 1.018365623020189
$ws.Cells.Item(9, 9).Value = 1.031581292705645
$ws.Cells.Item(9, 10).Value = 1.015118789941926
$ws.Cells.Item(9, 11).Value = 1.034802424814877
$ws.Cells.Item(9, 12).Value = 1.014586628014222
$ws.Cells.Item(9, 13).Value = 1.021586820279749
$ws.Cells.Item(9, 14).Value = 1.016560375059396

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.002174007647134
$ws.Cells.Item(10, 4).Value = 1.030358193332946
$ws.Cells.Item(10, 5).Value = 1.005255391497313
$ws.Cells.Item(10, 6).Value = 1.011674696245521
$ws.Cells.Item(10, 9).Value = 1.030852693276357
$ws.Cells.Item(10, 10).Value = 1.008748410132611
$ws.Cells.Item(10, 11).Value = 1.033824071269007
$ws.Cells.Item(10, 12).Value = 1.008814557471636
$ws.Cells.Item(10, 13).Value = 1.015209518669031
$ws.Cells.Item(10, 14).Value = 1.010180948580059

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 0.9990936587716028
$ws.Cells.Item(11, 4).Value = 1.029798032491476
$ws.Cells.Item(11, 5).Value = 1.00256621685782
$ws.Cells.Item(11, 6).Value = 1.008704654986604
$ws.Cells.Item(11, 9).Value = 1.030524600704798
$ws.Cells.Item(11, 10).Value = 1.005931024164105
$ws.Cells.Item(11, 11).Value = 1.033388281081302
$ws.Cells.Item(11, 12).Value = 1.006261857414141
$ws.Cells.Item(11, 13).Value = 1.012375986718477
$ws.Cells.Item(11, 14).Value = 1.007359561600319

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 0.9979397605317367
$ws.Cells.Item(12, 4).Value = 1.029588195173215
$ws.Cells.Item(12, 5).Value = 1.001558813250288
$ws.Cells.Item(12, 6).Value = 1.00758999962792
$ws.Cells.Item(12, 9).Value = 1.030400785558786
$ws.Cells.Item(12, 10).Value = 1.004875199893915
$ws.Cells.Item(12, 11).Value = 1.03322452840547
$ws.Cells.Item(12, 12).Value = 1.005305236422759
$ws.Cells.Item(12, 13).Value = 1.011312158993767
$ws.Cells.Item(12, 14).Value = 1.006302237938561

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 0.998187723763616
$ws.Cells.Item(13, 4).Value = 1.029633286953088
$ws.Cells.Item(13, 5).Value = 1.001775297702089
$ws.Cells.Item(13, 6).Value = 1.007829623626836
$ws.Cells.Item(13, 9).Value = 1.03042743332659
$ws.Cells.Item(13, 10).Value = 1.00510210733953
$ws.Cells.Item(13, 11).Value = 1.033259739933285
$ws.Cells.Item(13, 12).Value = 1.005510823551194
$ws.Cells.Item(13, 13).Value = 1.011540874509949
$ws.Cells.Item(13, 14).Value = 1.006529467618775

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 0.9989984784386189
$ws.Cells.Item(14, 4).Value = 1.029780723685592
$ws.Cells.Item(14, 5).Value = 1.002483121041648
$ws.Cells.Item(14, 6).Value = 1.008612753850853
$ws.Cells.Item(14, 9).Value = 1.030514406144183
$ws.Cells.Item(14, 10).Value = 1.005843942330964
$ws.Cells.Item(14, 11).Value = 1.033374783928876
$ws.Cells.Item(14, 12).Value = 1.006182957395123
$ws.Cells.Item(14, 13).Value = 1.01228828442813
$ws.Cells.Item(14, 14).Value = 1.007272356100987

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 0.9994967078469681
$ws.Cells.Item(15, 4).Value = 1.02987132824512
$ws.Cells.Item(15, 5).Value = 1.002918091459799
$ws.Cells.Item(15, 6).Value = 1.00909373327786
$ws.Cells.Item(15, 9).Value = 1.030567733323051
$ws.Cells.Item(15, 10).Value = 1.006299761844885
$ws.Cells.Item(15, 11).Value = 1.033445415420255
$ws.Cells.Item(15, 12).Value = 1.006595950640346
$ws.Cells.Item(15, 13).Value = 1.012747272046497
$ws.Cells.Item(15, 14).Value = 1.007728822930901

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.002377105190726
$ws.Cells.Item(16, 4).Value = 1.030395124453193
$ws.Cells.Item(16, 5).Value = 1.005432691983264
$ws.Cells.Item(16, 6).Value = 1.011870231048913
$ws.Cells.Item(16, 9).Value = 1.03087419767178
$ws.Cells.Item(16, 10).Value = 1.008934109521091
$ws.Cells.Item(16, 11).Value = 1.033852732539521
$ws.Cells.Item(16, 12).Value = 1.008982812477802
$ws.Cells.Item(16, 13).Value = 1.015396010026079
$ws.Cells.Item(16, 14).Value = 1.010366911682971

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.004167126169284
$ws.Cells.Item(17, 4).Value = 1.030720596852795
$ws.Cells.Item(17, 5).Value = 1.006995315812413
$ws.Cells.Item(17, 6).Value = 1.01359200792516
$ws.Cells.Item(17, 9).Value = 1.031063023601321
$ws.Cells.Item(17, 10).Value = 1.010570460760525
$ws.Cells.Item(17, 11).Value = 1.034104939832208
$ws.Cells.Item(17, 12).Value = 1.010465455434491
$ws.Cells.Item(17, 13).Value = 1.017037846795745
$ws.Cells.Item(17, 14).Value = 1.012005586728857

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.005205327480394
$ws.Cells.Item(18, 4).Value = 1.030909342989097
$ws.Cells.Item(18, 5).Value = 1.007901600747505
$ws.Cells.Item(18, 6).Value = 1.014589300579567
$ws.Cells.Item(18, 9).Value = 1.03117195023187
$ws.Cells.Item(18, 10).Value = 1.011519259701929
$ws.Cells.Item(18, 11).Value = 1.034250878999997
$ws.Cells.Item(18, 12).Value = 1.011325137517374
$ws.Cells.Item(18, 13).Value = 1.017988578900417
$ws.Cells.Item(18, 14).Value = 1.012955733073587

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.005558343262092
$ws.Cells.Item(19, 4).Value = 1.030973516194156
$ws.Cells.Item(19, 5).Value = 1.008209756734078
$ws.Cells.Item(19, 6).Value = 1.014928180796048
$ws.Cells.Item(19, 9).Value = 1.031208887326135
$ws.Cells.Item(19, 10).Value = 1.011841829669849
$ws.Cells.Item(19, 11).Value = 1.034300443782313
$ws.Cells.Item(19, 12).Value = 1.011617411009766
$ws.Cells.Item(19, 13).Value = 1.018311594102305
$ws.Cells.Item(19, 14).Value = 1.013278761127861

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.003975686343305
$ws.Cells.Item(20, 4).Value = 1.030685790539495
$ws.Cells.Item(20, 5).Value = 1.006828198556951
$ws.Cells.Item(20, 6).Value = 1.01340800421455
$ws.Cells.Item(20, 9).Value = 1.031042890138423
$ws.Cells.Item(20, 10).Value = 1.010395484146244
$ws.Cells.Item(20, 11).Value = 1.034078001681863
$ws.Cells.Item(20, 12).Value = 1.010306914248999
$ws.Cells.Item(20, 13).Value = 1.016862412844534
$ws.Cells.Item(20, 14).Value = 1.011830361627714

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 0.9987600038442006
$ws.Cells.Item(21, 4).Value = 1.0297373565282
$ws.Cells.Item(21, 5).Value = 1.002274923619514
$ws.Cells.Item(21, 6).Value = 1.008382461832021
$ws.Cells.Item(21, 9).Value = 1.030488849008523
$ws.Cells.Item(21, 10).Value = 1.005625751592298
$ws.Cells.Item(21, 11).Value = 1.033340958713663
$ws.Cells.Item(21, 12).Value = 1.005985267052826
$ws.Cells.Item(21, 13).Value = 1.012068507650695
$ws.Cells.Item(21, 14).Value = 1.007053855506444

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 0.9954241946905454
$ws.Cells.Item(22, 4).Value = 1.029130775539129
$ws.Cells.Item(22, 5).Value = 0.9993625465572125
$ws.Cells.Item(22, 6).Value = 1.005156183717003
$ws.Cells.Item(22, 9).Value = 1.030129207873483
$ws.Cells.Item(22, 10).Value = 1.002572649467712
$ws.Cells.Item(22, 11).Value = 1.032866639251936
$ws.Cells.Item(22, 12).Value = 1.003219050401886
$ws.Cells.Item(22, 13).Value = 1.008988583318819
$ws.Cells.Item(22, 14).Value = 1.003996417626646

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 0.9971980999851382
$ws.Cells.Item(23, 4).Value = 1.029453327828638
$ws.Cells.Item(23, 5).Value = 1.000911300605785
$ws.Cells.Item(23, 6).Value = 1.006872978289883
$ws.Cells.Item(23, 9).Value = 1.030320949667178
$ws.Cells.Item(23, 10).Value = 1.004196454236201
$ws.Cells.Item(23, 11).Value = 1.03311913834926
$ws.Cells.Item(23, 12).Value = 1.004690267665705
$ws.Cells.Item(23, 13).Value = 1.010627719795297
$ws.Cells.Item(23, 14).Value = 1.00562252838416

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.004062207919258
$ws.Cells.Item(24, 4).Value = 1.030701521398709
$ws.Cells.Item(24, 5).Value = 1.006903727588298
$ws.Cells.Item(24, 6).Value = 1.013491169137747
$ws.Cells.Item(24, 9).Value = 1.031051991326338
$ws.Cells.Item(24, 10).Value = 1.010474565995346
$ws.Cells.Item(24, 11).Value = 1.034090177481454
$ws.Cells.Item(24, 12).Value = 1.010378567974035
$ws.Cells.Item(24, 13).Value = 1.016941705255876
$ws.Cells.Item(24, 14).Value = 1.011909555782113

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.01179818615028
$ws.Cells.Item(25, 4).Value = 1.032107137287027
$ws.Cells.Item(25, 5).Value = 1.013656232607152
$ws.Cells.Item(25, 6).Value = 1.020898532347099
$ws.Cells.Item(25, 9).Value = 1.031852927474522
$ws.Cells.Item(25, 10).Value = 1.017539466192214
$ws.Cells.Item(25, 11).Value = 1.035171317594516
$ws.Cells.Item(25, 12).Value = 1.016780019264983
$ws.Cells.Item(25, 13).Value = 1.023998641101585
$ws.Cells.Item(25, 14).Value = 1.018984488947615
